$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text cell updates (column B)
$ws.Range("B2").Value = "<are>"
$ws.Range("B6").Value = "<sevence>"
$ws.Range("B12").Value = "<him>"

# Numeric cell updates (column C)
$ws.Range("C3").Value = 22
$ws.Range("C4").Value = 13
$ws.Range("C5").Value = 18
$ws.Range("C6").Value = 17
$ws.Range("C7").Value = 9
$ws.Range("C8").Value = 18
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 12
$ws.Range("C13").Value = 16
$ws.Range("C14").Value = 8
$ws.Range("C15").Value = 8
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 15
